$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.684.14"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.632.17"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "213.36"
$ws.Range("E5").Value = "  +0.57%  "
Set-TextValue $ws.Range("D6") "0.501"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("E9").Value = "  +0.95%  "
Set-TextValue $ws.Range("D10") "19.24"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "1.859.30"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.627.29"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "26.671.40"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("E18").Value = "  +2.13%  "
Set-TextValue $ws.Range("D19") "218.96"
$ws.Range("E19").Value = "  +8.29%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.95%  "
Set-TextValue $ws.Range("D23") "9.37"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("E24").Value = "  +4.54%  "
Set-TextValue $ws.Range("D25") "147.52"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("E28").Value = "  +4.39%  "
Set-TextValue $ws.Range("D29") "15.58"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "1.216.48"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  -1.91%  "
Set-TextValue $ws.Range("D42") "0.795"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "1.770.42"
$ws.Range("E44").Value = "  +0.63%  "
Set-TextValue $ws.Range("D45") "92.78"
$ws.Range("E45").Value = "  +0.41%  "
Set-TextValue $ws.Range("D46") "1.56"
$ws.Range("E46").Value = "  +2.57%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -0.94%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D48") "55.07"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("E49").Value = "  +0.61%  "
Set-TextValue $ws.Range("D50") "7.60"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("E51").Value = "  -0.08%  "
